$wb = $excel.ActiveWorkbook

# Bus sheet: update computed bus voltage/angle results (v0, a0) and move the selection
$busSheet = $wb.Worksheets.Item("Bus")
$busSheet.Activate()
$busSheet.Range("I2").Value = 0.57025491716260168
$busSheet.Range("H3").Value = 0.99761
$busSheet.Range("I3").Value = 0.36874618304434781
$busSheet.Range("H3:I3").Select()

# PQ sheet: correct p0/q0/vmax units (480 -> 0.48 pu, q0 results reset to 0) and move the selection
$pqSheet = $wb.Worksheets.Item("PQ")
$pqSheet.Activate()
$pqSheet.Range("F2").Value = 0.48
$pqSheet.Range("G2").Value = 0
$pqSheet.Range("H2").Value = 0
$pqSheet.Range("H3").Select()

# Slack sheet: no data changes, only move the selection
$slackSheet = $wb.Worksheets.Item("Slack")
$slackSheet.Activate()
$slackSheet.Range("A2").Select()

# Line sheet: fix r/x (x,b) magnitudes by a factor of 10 and move the selection
$lineSheet = $wb.Worksheets.Item("Line")
$lineSheet.Activate()
$lineSheet.Range("L2").Value = 0.005
$lineSheet.Range("M2").Value = 0.0075
$lineSheet.Range("M2").Select()

# Synchronverter sheet: no data changes, only move the selection (keeps this tab active, as before)
$synchSheet = $wb.Worksheets.Item("Synchronverter")
$synchSheet.Activate()
$synchSheet.Range("G31").Select()
